$d = $word.ActiveDocument

# 1. Merge the 4 runs of "Master of Science (MSc) In Pure Mathematics (Gold " / "Medallist" /
#    ") University of Calcutta" / "." into a single run.
$d.Content.Find.Execute("Master of Science (MSc) In Pure Mathematics (Gold Medallist) University of Calcutta.", `
  $true, $false, $false, $false, $false, $true, 1, $false, `
  "Master of Science (MSc) In Pure Mathematics (Gold Medallist) University of Calcutta.", 2)

# 2. Merge the 6 runs of "Stood First Class First " / "in" / " " / "BSc (" / "Mathematics Major)" /
#    "." into a single run.
$d.Content.Find.Execute("Stood First Class First in BSc (Mathematics Major).", `
  $true, $false, $false, $false, $false, $true, 1, $false, `
  "Stood First Class First in BSc (Mathematics Major).", 2)

# 3. Merge the 2 runs "MSc" / " (Pure Mathematics)  " into a single run.
$d.Content.Find.Execute("MSc (Pure Mathematics)  ", `
  $true, $false, $false, $false, $false, $true, 1, $false, `
  "MSc (Pure Mathematics)  ", 2)

# 4. Merge the 2 runs "Professional " / "Experience " into a single run (paragraph/structure
#    otherwise unchanged).
$d.Content.Find.Execute("Professional Experience ", `
  $true, $false, $false, $false, $false, $true, 1, $false, `
  "Professional Experience ", 2)

# 5. Merge the 3 runs "students in subjects like " / "Mathematics, Physics & Chemistry" / " for"
#    into a single run, *without* pulling in the preceding separate "Mentoring " run. Replacing
#    just the middle run's own text (identical content) causes the engine to coalesce it with its
#    immediate identically-formatted neighbours on both sides, stopping at "Mentoring ".
$d.Content.Find.Execute("Mathematics, Physics & Chemistry", `
  $true, $false, $false, $false, $false, $true, 1, $false, `
  "Mathematics, Physics & Chemistry", 2)
